$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Diebold-Mariano statistics (column C) and P-values (column D)
$ws.Range("C2").Value = 0.544108002761951
$ws.Range("D2").Value = 0.5899178377614889

$ws.Range("C3").Value = -0.1717742808392656
$ws.Range("D3").Value = 0.8646333209788317

$ws.Range("C4").Value = -1.125771030857224
$ws.Range("D4").Value = 0.2681464611039448

$ws.Range("C5").Value = -0.2711000261073233
$ws.Range("D5").Value = 0.7879532771660034

$ws.Range("C6").Value = -0.7046827443844143
$ws.Range("D6").Value = 0.4858070976686637

$ws.Range("C7").Value = -1.101916103017213
$ws.Range("D7").Value = 0.2782376315500676

$ws.Range("C8").Value = -0.6411013706959228
$ws.Range("D8").Value = 0.5257550380154998

$ws.Range("C9").Value = -1.197692430595533
$ws.Range("D9").Value = 0.239320000207796

$ws.Range("C10").Value = -0.2234485918706526
$ws.Range("D10").Value = 0.8245235438414715

$ws.Range("C11").Value = 1.287159237078364
$ws.Range("D11").Value = 0.2067327862896275
